$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pub runs")

# 1. Insert a new column before column E (shifts old E:R -> F:S)
$ws.Columns("E:E").Insert()

# 2. Insert a new row before row 69 (shifts old row69/70/73 -> 70/71/74)
$ws.Rows("69:69").Insert()

# 3. Fill in the new pub-run entry's strings + the "Km" header, in the exact
#    write order needed so new shared-string indices land where the target
#    workbook expects them (188=The Hoptimist, 189=Spondon, 190=Km,
#    191=S4007/TP2628 - Crow Wood Farm, 192=Attacked by herd of deer)
$ws.Range("B69").Value = "The Hoptimist"
$ws.Range("C69").Value = "Spondon"
$ws.Range("E1").Value = "Km"
$ws.Range("S69").Value = "S4007/TP2628 - Crow Wood Farm"
$ws.Range("P69").Value = "Attacked by herd of deer"

$ws.Range("A69").Value = 44958
$ws.Range("D69").Value = "start/end at pub"
$ws.Range("E69").Value = 5.01
$ws.Range("F69").Formula = "=E69*0.6213712"
$ws.Range("G69").Value = 0.028692129629629633
$ws.Range("H69").Formula = "=G69/F69"
$ws.Range("I69").Value = 1
$ws.Range("O69").Value = 1
$ws.Range("Q69").Formula = "=SUM(I69:O69)*F69"

# 5. Update the two grand-total formulas that were hand-extended to include
#    the new row 69 (CL and RM counts)
$ws.Range("I71").Formula = "=SUM(I3:I69)"
$ws.Range("O71").Formula = "=SUM(O3:O69)"
